$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.600.32"
$ws.Range("E2").Value = "  -0.63%  "

$ws.Range("D3").Value = "3.442.30"
$ws.Range("E3").Value = "  -2.33%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.94%  "

$ws.Range("D7").Value = "3.440.97"
$ws.Range("E7").Value = "  -2.39%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.86%  "

$ws.Range("E11").Value = "  -8.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.378"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.59%  "

$ws.Range("D13").Value = "4.028.56"
$ws.Range("E13").Value = "  -2.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000181"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.01%  "

$ws.Range("D16").Value = "3.437.59"
$ws.Range("E16").Value = "  -2.28%  "

$ws.Range("D17").Value = "65.576.08"
$ws.Range("E17").Value = "  -0.66%  "

$ws.Range("E18").Value = "  -1.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -9.52%  "

$ws.Range("E20").Value = "  -4.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "393.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.556"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.68%  "

$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").Value = "3.579.67"
$ws.Range("E26").Value = "  -2.48%  "

$ws.Range("E27").Value = "  -6.33%  "

$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("E29").Value = "  -5.91%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.59%  "

$ws.Range("E31").Value = "  -8.64%  "

$ws.Range("D32").Value = "3.449.01"
$ws.Range("E32").Value = "  -2.19%  "

$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("E34").Value = "  -4.74%  "

$ws.Range("E35").Value = "  -4.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "172.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.84%  "

$ws.Range("E38").Value = "  -5.61%  "

$ws.Range("E39").Value = "  -6.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0765"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.825"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.84%  "

$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("E45").Value = "  -12.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.09%  "

$ws.Range("E50").Value = "  -11.84%  "

$ws.Range("D51").Value = "2.215.24"
$ws.Range("E51").Value = "  -6.24%  "
